$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A5").Value = "bh"
$ws.Range("A6").Value = "Piyali"
